# Generate Report for Handoff
#
# This script takes the localization-status workbook from a state where
# only one source file (0d1f9c3b-...md) has been handed off, to a state
# where a second handoff round has produced two additional dependent
# files (two .png images) alongside the original .md, refreshing the
# handoff timestamps/file hashes along the way.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Identifiers used throughout
# ---------------------------------------------------------------------
$mdFile      = "6c65d407-12d1-47e5-9b53-fc290c4112cf.md"
$pngFile1    = "967f602c-39ee-420a-adf0-dc15cf2fcf9c.png"
$pngFile2    = "caa9e922-df7c-43ba-a579-9697b9ca4822.png"

$zhXlf       = "6c65d407-12d1-47e5-9b53-fc290c4112cf.c0a6579dbc28ebdefd6a05e730dadf21335a4523.zh-cn.xlf"
$deXlf       = "6c65d407-12d1-47e5-9b53-fc290c4112cf.c0a6579dbc28ebdefd6a05e730dadf21335a4523.de-de.xlf"
$pngHash1    = "22bc2beee41db635d56058ec1fc5656ac11755df.png"
$pngHash2    = "c68a3dc8e662f9d1e81af1f4ba7728a0d53de193.png"

$handoffDate      = "2016-03-22 09:01:37"
$handoffDatetime  = "2016-03-22 09:01:27"
$zeroDate         = "0001-01-01 00:00:00"
$readyStatus      = "Ready for handoff"

$dependencyFrom   = "e2e\" + $mdFile

$mdUrlBase  = "https://github.com/OpenLocalizationTest/oltest/blob/f588d4f5f45fd06c0b006db141dd9f5c2136af93/e2e/"
$zhUrlBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed845000f641fdfadfd7f9b75f9fa092026ad133/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deUrlBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0ab79a0b60388d57f51dbfbe372c42b0ea2ecf7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================

# -- refresh row 2 (existing handoff, new round) --
$ws1.Range("A2").Value2 = $mdFile
$ws1.Range("B2").Value2 = $readyStatus
$ws1.Range("C2").Value2 = $readyStatus
$ws1.Range("D2").Value2 = $handoffDate

# -- new row 3 / row 4 for the two dependent png files --
$ws1.Range("A3").Value2 = $pngFile1
$ws1.Range("B3").Value2 = $readyStatus
$ws1.Range("C3").Value2 = $readyStatus
$ws1.Range("D3").Value2 = $handoffDate
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("A4").Value2 = $pngFile2
$ws1.Range("B4").Value2 = $readyStatus
$ws1.Range("C4").Value2 = $readyStatus
$ws1.Range("D4").Value2 = $handoffDate
$ws1.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- rebuild the hyperlinks in order: A2, A3, A4 --
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdUrlBase + $mdFile), "", "", $mdFile)
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdUrlBase + $pngFile1), "", "", $pngFile1)
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($mdUrlBase + $pngFile2), "", "", $pngFile2)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================

# -- refresh row 2 --
$ws2.Range("A2").Value2 = $mdFile
$ws2.Range("B2").Value2 = ".md"
$ws2.Range("C2").Value2 = $readyStatus
$ws2.Range("D2").Value2 = $zhXlf
$ws2.Range("E2").Value2 = $handoffDatetime
$ws2.Range("H2").Value2 = $zeroDate
$ws2.Range("J2").Value2 = "Include"

# -- new row 3 (png dependency 1) --
$ws2.Range("A3").Value2 = $pngFile1
$ws2.Range("B3").Value2 = ".png"
$ws2.Range("C3").Value2 = $readyStatus
$ws2.Range("D3").Value2 = $pngHash1
$ws2.Range("E3").Value2 = $handoffDatetime
$ws2.Range("H3").Value2 = $zeroDate
$ws2.Range("J3").Value2 = "IsDependency"
$ws2.Range("K3").Value2 = $dependencyFrom

# -- new row 4 (png dependency 2) --
$ws2.Range("A4").Value2 = $pngFile2
$ws2.Range("B4").Value2 = ".png"
$ws2.Range("C4").Value2 = $readyStatus
$ws2.Range("D4").Value2 = $pngHash2
$ws2.Range("E4").Value2 = $handoffDatetime
$ws2.Range("H4").Value2 = $zeroDate
$ws2.Range("J4").Value2 = "IsDependency"
$ws2.Range("K4").Value2 = $dependencyFrom

# -- rebuild hyperlinks in order: A2, D2, A3, D3, A4, D4 --
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdUrlBase + $mdFile), "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("D2"), ($zhUrlBase + $zhXlf), "", "", $zhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdUrlBase + $pngFile1), "", "", $pngFile1)
$ws2.Hyperlinks.Add($ws2.Range("D3"), ($zhUrlBase + $pngHash1), "", "", $pngHash1)
$ws2.Hyperlinks.Add($ws2.Range("A4"), ($mdUrlBase + $pngFile2), "", "", $pngFile2)
$ws2.Hyperlinks.Add($ws2.Range("D4"), ($zhUrlBase + $pngHash2), "", "", $pngHash2)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================

# -- refresh row 2 --
$ws3.Range("A2").Value2 = $mdFile
$ws3.Range("B2").Value2 = ".md"
$ws3.Range("C2").Value2 = $readyStatus
$ws3.Range("D2").Value2 = $deXlf
$ws3.Range("E2").Value2 = "de-de"
$ws3.Range("H2").Value2 = $zeroDate
$ws3.Range("J2").Value2 = "Include"

# -- new row 3 (png dependency 1) --
$ws3.Range("A3").Value2 = $pngFile1
$ws3.Range("B3").Value2 = ".png"
$ws3.Range("C3").Value2 = $readyStatus
$ws3.Range("D3").Value2 = $pngHash1
$ws3.Range("E3").Value2 = "de-de"
$ws3.Range("H3").Value2 = $zeroDate
$ws3.Range("J3").Value2 = "IsDependency"
$ws3.Range("K3").Value2 = $dependencyFrom

# -- new row 4 (png dependency 2) --
$ws3.Range("A4").Value2 = $pngFile2
$ws3.Range("B4").Value2 = ".png"
$ws3.Range("C4").Value2 = $readyStatus
$ws3.Range("D4").Value2 = $pngHash2
$ws3.Range("E4").Value2 = "de-de"
$ws3.Range("H4").Value2 = $zeroDate
$ws3.Range("J4").Value2 = "IsDependency"
$ws3.Range("K4").Value2 = $dependencyFrom

# -- rebuild hyperlinks in order: A2, D2, A3, D3, A4, D4 --
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdUrlBase + $mdFile), "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("D2"), ($deUrlBase + $deXlf), "", "", $deXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdUrlBase + $pngFile1), "", "", $pngFile1)
$ws3.Hyperlinks.Add($ws3.Range("D3"), ($deUrlBase + $pngHash1), "", "", $pngHash1)
$ws3.Hyperlinks.Add($ws3.Range("A4"), ($mdUrlBase + $pngFile2), "", "", $pngFile2)
$ws3.Hyperlinks.Add($ws3.Range("D4"), ($deUrlBase + $pngHash2), "", "", $pngHash2)

Write-Host "Report generated for handoff"
